$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# --- Update existing data row (row 2) with the new error codes ---
$t.Rows.Item(2).Cells.Item(2).Range.Text = "INVALID USERNAME (2001)"
$t.Rows.Item(2).Cells.Item(3).Range.Text = "USER_EXISTED (3001)"
$t.Rows.Item(2).Cells.Item(4).Range.Text = "USER NOT FOUND (4001)"

# --- Add a new row for the 2002 / 4002 error codes ---
$row3 = $t.Rows.Add()
$row3.Cells.Item(2).Range.Text = "INVALID PASSWORD (2002)"
$row3.Cells.Item(4).Range.Text = "ROLE NOT FOUND (4002)"

# --- Add a new row for the 2003 error code ---
$row4 = $t.Rows.Add()
$row4.Cells.Item(2).Range.Text = "REQUIRED EMAIL (2003)"

# --- Resize the table columns to match the updated layout ---
$t.Columns.Item(1).Width = 97.55
$t.Columns.Item(2).Width = 72.7
$t.Columns.Item(3).Width = 88.85
$t.Columns.Item(4).Width = 96.7
$t.Columns.Item(5).Width = 112.85
$t.Columns.Item(6).Width = 112.55
